$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.270.22'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.790.90'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.68%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.008'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '307.56'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.75%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4547'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.68%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07065'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8693'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07767'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.29'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.784.93'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.259'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.302'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '84.54'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -7.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.010'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008494'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.008'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.343.10'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.08'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.962'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.017.09'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -5.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.47'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.983'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.27'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.75'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.024'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.08'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.823'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08646'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.034'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.429'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7110'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -9.07%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.097'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.32%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.621'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.008'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.078'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01933'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05081'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.865'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.870'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4908'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1516'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.943'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.16%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4557'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.798'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '99.67'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.576'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05936'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.25%  '
